# fix(docx): fix OOXMLValidator error on KeywordTok output
#
# wml.xsd's CT_RPr requires the boolean toggle properties (<w:b/>,
# <w:i/>, ...) to precede <w:color/> inside <w:rPr>. Several of the
# Pandoc "highlighting" character styles in styles.xml had them in the
# wrong order (<w:color/> before <w:b/>/<w:i/>), which xmllint let
# through but OOXMLValidatorCLI flagged as
# Sch_UnexpectedElementContentExpectingComplex.
#
# Re-apply the existing bold/italic formatting on each affected style
# so its <w:rPr> gets re-serialized with elements in schema order.

$d = $word.ActiveDocument
$styles = $d.Styles

# <w:color/><w:b/>            -> <w:b/><w:color/>
$styles.Item("KeywordTok").Font.Bold = 1
$styles.Item("ImportTok").Font.Bold = 1
$styles.Item("ControlFlowTok").Font.Bold = 1
$styles.Item("AlertTok").Font.Bold = 1
$styles.Item("ErrorTok").Font.Bold = 1

# <w:color/><w:i/>            -> <w:i/><w:color/>
$styles.Item("CommentTok").Font.Italic = 1
$styles.Item("DocumentationTok").Font.Italic = 1

# <w:color/><w:b/><w:i/>      -> <w:b/><w:i/><w:color/>
$styles.Item("AnnotationTok").Font.Bold = 1
$styles.Item("AnnotationTok").Font.Italic = 1
$styles.Item("CommentVarTok").Font.Bold = 1
$styles.Item("CommentVarTok").Font.Italic = 1
$styles.Item("InformationTok").Font.Bold = 1
$styles.Item("InformationTok").Font.Italic = 1
$styles.Item("WarningTok").Font.Bold = 1
$styles.Item("WarningTok").Font.Italic = 1
